# Update "想去人数" (number of people interested) values in column F
# for the 展览 (Exhibition) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - row => new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4294
$ws1.Range("F3").Value = 2433
$ws1.Range("F5").Value = 20
$ws1.Range("F7").Value = 56
$ws1.Range("F11").Value = 154
$ws1.Range("F12").Value = 1595
$ws1.Range("F14").Value = 3337
$ws1.Range("F15").Value = 226

# Sheet "全部类型" (All types) - row => new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4294
$ws4.Range("F3").Value = 2433
$ws4.Range("F5").Value = 20
$ws4.Range("F8").Value = 56
$ws4.Range("F13").Value = 154
$ws4.Range("F16").Value = 1595
$ws4.Range("F18").Value = 3337
$ws4.Range("F19").Value = 226
